# Re-apply the (built-in) "Medium Style 2 - Accent 1" table style
# ({96DCA32E-0A28-41AC-A9E3-E4AC378BD2C5}) to every table in the deck.
# Previously every table carried the deck's single custom table style
# ({18D845A6-A52B-4B88-AA2D-ED0B2E175017}) defined in tableStyles.xml.

$newStyleId = "{96DCA32E-0A28-41AC-A9E3-E4AC378BD2C5}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $true)
        }
    }
}
